$d = $word.ActiveDocument

# The work log is the first table in the document.
$t = $d.Tables.Item(1)

# Add a new row for 3/4/24
$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "3/4/24"
$row1.Cells.Item(2).Range.Text = "Investigating mapper service and sending data to the API. Nothing worth committing."
$row1.Cells.Item(3).Range.Text = "4"

# Add a new row for 4/4/24
$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "4/4/24"
$row2.Cells.Item(2).Range.Text = "Successfully sending data to the mapper service. Made adjustments to the UI dashboard for a blinking LED as pH warning. Saved copy of Node-red flows to src folder."
$row2.Cells.Item(3).Range.Text = "4"
